# Apply delete test case changes to the RestAssured sheet (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RestAssured")

# Update isbn (C) and aisle (D) values for rows 2-5
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 2

$ws.Range("C3").Value = 14
$ws.Range("D3").Value = 2

$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 2

$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 2

# Update the active cell selection from D14 to E14
$ws.Range("E14").Select()
